# Applies the commit's corrections to the SSU "lost password recovery"
# page specification.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Grammar fix: "opisuju" -> "opisuje"
#    "U ovom odeljku se opisuju glavni uspešni scenario ..." becomes
#    "U ovom odeljku se opisuje glavni uspešni scenario ..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "U ovom odeljku se opisuju glavni uspešni scenario",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "U ovom odeljku se opisuje glavni uspešni scenario", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Drop the "i podnožje (footer tag) stranice" clause - the page is
#    now described as only having a header tag.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    " (header tag) i podnožje (footer tag) stranice; polj",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " (header tag); polj", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) The button caption changes from the Serbian "Pošalji" to the
#    English "Send" everywhere it is mentioned.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Pošalji",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Send", 2) | Out-Null

$d.Content.Find.Execute(
    "Pošalji",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Send", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Minor run clean-up: ", u " is written as one contiguous run in the
#    "Sistem šalje poruku ..." sentence (paragraph 110 - restrict the
#    search to this paragraph since the same phrase also occurs, verbatim,
#    in an earlier summary paragraph that must stay untouched).
# ---------------------------------------------------------------------
$sistemPar = $d.Paragraphs(110)
$sistemPar.Range.Find.Execute(
    ", u ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", u ", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Append the missing "consequences" paragraph at the very end of the
#    document body (after the "Posledice" heading).
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$lastPar = $d.Paragraphs.Last
$lastPar.Range.Text = "Korisniku je poslat mejl sa linkom i instrukcijama za promenu šifre."
$lastPar.Style = "PSI Normal"
